# Update the Hellwig ranking values (B2:B17) with the corrected method
# results, apply wrap-text formatting to that range (creates the new
# cellXfs entry used by those cells), resize column B slightly, and move
# the active selection to H6 - mirroring the authored commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    0.243164,
    0.246206,
    0.30354,
    0.208139,
    0.339088,
    0.667013,
    0.795715,
    0.140076,
    0.225525,
    0.561972,
    0.408323,
    0.421173,
    0.339127,
    0.112509,
    0.36864,
    0.3473
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}

# New data rows wrap their text (introduces the 4th cellXfs entry, s="3").
$ws.Range("B2:B17").WrapText = $true

# Column B widened slightly to fit the new formatting.
$ws.Columns.Item(2).ColumnWidth = 7.09

# Leave the selection on H6, as in the saved workbook.
$ws.Range("H6").Select()
